$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Chennai_SouthAfrica: Matches 2 -> 3, Lost 2 -> 3
$ws.Range("B3").Value = 3.0
$ws.Range("D3").Value = 3.0

# Row 5 - RoyalChallengers_Bangladesh: Matches 3 -> 4, Won 1 -> 2, Points 2 -> 4
$ws.Range("B5").Value = 4.0
$ws.Range("C5").Value = 2.0
$ws.Range("E5").Value = 4.0

# Row 7 - Rajastan_Australia: Matches 4 -> 5, Won 2 -> 3, Points 4 -> 6
$ws.Range("B7").Value = 5.0
$ws.Range("C7").Value = 3.0
$ws.Range("E7").Value = 6.0

# Row 9 - Punjab_Pakistan: Matches 7 -> 8, Lost 3 -> 4
$ws.Range("B9").Value = 8.0
$ws.Range("D9").Value = 4.0
